$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the existing table to cover the new range (adds 2 columns and 15 rows)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J22"))

# The source workbook stores a handful of intentionally-blank table cells as
# typed-but-empty shared-string cells. Re-clear them so they stay genuinely
# blank (and are not misread as pointing at shared-string index 0).
$ws.Range("E2").ClearContents()
$ws.Range("A3:D3").ClearContents()
$ws.Range("A4:D4").ClearContents()
$ws.Range("A5:E5").ClearContents()
$ws.Range("A6:E6").ClearContents()
$ws.Range("A7:E7").ClearContents()

# Set header cells (new columns); this also renames the corresponding ListColumns
$ws.Range("I1").Value = 'Approved/Rejected'
$ws.Range("J1").Value = 'ReasonToReject'

# Populate data rows 2-22 for columns A-H with the new/updated test-case content
$ws.Range("A2").Value = 'TestScenario_1'
$ws.Range("B2").Value = 'TestScenario_1.TestCase_1'
$ws.Range("C2").Value = 'New Account'
$ws.Range("D2").Value = 'User Needs to Login to Salesforce, from the browser with correct credentials'
$ws.Range("F2").Value = 'Step 1'
$ws.Range("G2").Value = 'Click Account tab, and click on New button'
$ws.Range("H2").Value = 'Shows fields to enter to create a new Account'
$ws.Range("E3").Value = 'Valid value for required field Name'
$ws.Range("F3").Value = 'Step 2'
$ws.Range("G3").Value = 'Enter valid value for Name'
$ws.Range("H3").Value = 'Value accepted for Name'
$ws.Range("E4").Value = 'Valid value for required field AnnualRevenue'
$ws.Range("F4").Value = 'Step 3'
$ws.Range("G4").Value = 'Enter valid value for AnnualRevenue, value should be greater or equal 50000.45'
$ws.Range("H4").Value = 'Value accepted for AnnualRevenue'
$ws.Range("F5").Value = 'Step 4'
$ws.Range("G5").Value = 'Click on Save button to save Account with fields'
$ws.Range("H5").Value = 'Saved changes successfully for the Account'
$ws.Range("F6").Value = 'Step 5'
$ws.Range("G6").Value = 'Click on ''Submit for Approval'' button to submit the record for Approval.'
$ws.Range("H6").Value = 'Alert message box will be displayed for confirmation from Salesforce.'
$ws.Range("F7").Value = 'Step 6'
$ws.Range("G7").Value = 'Click on ''OK'' button to submit the record for Approval.'
$ws.Range("H7").Value = 'Unable to Submit for Approval message will be displayed if this record does not meet the entry criteria. Otherwise, this record will be displayed under Approval History section with the status ''Pending''.'
$ws.Range("A8").Value = 'TestScenario_2'
$ws.Range("B8").Value = 'TestScenario_2.TestCase_1'
$ws.Range("C8").Value = 'View Account'
$ws.Range("D8").Value = 'User Needs to Login to Salesforce, from the browser with correct credentials'
$ws.Range("F8").Value = 'Step 1'
$ws.Range("G8").Value = 'Click Account tab,  and click on existing Account name to view'
$ws.Range("H8").Value = 'Shows fields for Account selected'
$ws.Range("F9").Value = 'Step 2'
$ws.Range("G9").Value = 'Click on Account name to View the Details'
$ws.Range("H9").Value = 'Details of Account successfully displayed'
$ws.Range("A10").Value = 'TestScenario_3'
$ws.Range("B10").Value = 'TestScenario_3.TestCase_1'
$ws.Range("C10").Value = 'Edit Account'
$ws.Range("D10").Value = 'User Needs to Login to Salesforce, from the browser with correct credentials'
$ws.Range("F10").Value = 'Step 1'
$ws.Range("G10").Value = 'Click Account tab,  and click on existing Account to modify'
$ws.Range("H10").Value = 'Shows fields to modify for an existing Account'
$ws.Range("E11").Value = 'Valid value for required field Name'
$ws.Range("F11").Value = 'Step 2'
$ws.Range("G11").Value = 'Enter valid value for Name'
$ws.Range("H11").Value = 'Value accepted for Name'
$ws.Range("E12").Value = 'Valid value for required field AnnualRevenue'
$ws.Range("F12").Value = 'Step 3'
$ws.Range("G12").Value = 'Enter valid value for AnnualRevenue, value should be greater or equal 50000.45'
$ws.Range("H12").Value = 'Value accepted for AnnualRevenue'
$ws.Range("F13").Value = 'Step 4'
$ws.Range("G13").Value = 'Click on Save button to save Account with fields'
$ws.Range("H13").Value = 'Saved changes successfully for the Account'
$ws.Range("F14").Value = 'Step 5'
$ws.Range("G14").Value = 'Click on ''Submit for Approval'' button to submit the record for Approval.'
$ws.Range("H14").Value = 'Alert message box will be displayed for confirmation from Salesforce.'
$ws.Range("F15").Value = 'Step 6'
$ws.Range("G15").Value = 'Click on ''OK'' button to submit the record for Approval.'
$ws.Range("H15").Value = 'Unable to Submit for Approval message will be displayed if this record does not meet the entry criteria. Otherwise, this record will be displayed under Approval History section with the status ''Pending''.'
$ws.Range("A16").Value = 'TestScenario_4'
$ws.Range("B16").Value = 'TestScenario_4.TestCase_1'
$ws.Range("C16").Value = 'Delete Account'
$ws.Range("D16").Value = 'User Needs to Login to Salesforce, from the browser with correct credentials'
$ws.Range("F16").Value = 'Step 1'
$ws.Range("G16").Value = 'Click Account tab,  and click on existing Account to delete'
$ws.Range("H16").Value = 'Show popup to confirm deletion of Account'
$ws.Range("F17").Value = 'Step 2'
$ws.Range("G17").Value = 'Click yes on confirm dialog to Delete the Account'
$ws.Range("H17").Value = 'Deleted the Account successfully'
$ws.Range("A18").Value = 'TestScenario_1'
$ws.Range("B18").Value = 'TestScenario_1.TestCase_1'
$ws.Range("C18").Value = 'Approve/Reject Account'
$ws.Range("D18").Value = 'User Needs to Login to Salesforce, from the browser with correct credentials'
$ws.Range("F18").Value = 'Step 1'
$ws.Range("G18").Value = 'Click Account tab, '
$ws.Range("F19").Value = 'Step 2'
$ws.Range("G19").Value = 'Click on Account name to Approve/Reject'
$ws.Range("H19").Value = 'Details of Account successfully displayed'
$ws.Range("F20").Value = 'Step 3'
$ws.Range("G20").Value = 'Click on Approve/Reject link from ''Approval History Section'' to Approve/Reject Account'
$ws.Range("H20").Value = 'Details of ''Approve/Reject Approval Request'' successfully displayed'
$ws.Range("F21").Value = 'Step 4'
$ws.Range("G21").Value = 'Enter Comments to notify the user (if any)'
$ws.Range("H21").Value = 'Value accepted for Comments'
$ws.Range("F22").Value = 'Step 5'
$ws.Range("G22").Value = 'Click on Approve/Reject button to Approve/Reject Account'
$ws.Range("H22").Value = 'Overall status (Approved/Rejected) will be displayed under ''Approval History'' section'

# Adjust column widths to match the new layout
$ws.Columns.Item(3).ColumnWidth = 23.666666666666668
$ws.Columns.Item(7).ColumnWidth = 79
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(10).ColumnWidth = 16.166666666666668
